$d = $word.ActiveDocument
$sec = $d.Sections(1)

# --- First-page header (header1.xml): BTec logo image1.jpg -> image2.jpg ---
$hdrFirst = $sec.Headers(2)
if ($hdrFirst.Exists -and $hdrFirst.Range.InlineShapes.Count -ge 1) {
    $btecLogo = $hdrFirst.Range.InlineShapes(1)
    $btecLogo.Name = "image2.jpg"
}

# --- Default footer (footer2.xml): Pearson logo image2.png -> image1.png ---
$ftrDefault = $sec.Footers(1)
if ($ftrDefault.Exists -and $ftrDefault.Range.InlineShapes.Count -ge 1) {
    $pearsonLogoA = $ftrDefault.Range.InlineShapes(1)
    $pearsonLogoA.Name = "image1.png"
}

# --- First-page footer (footer1.xml): Pearson logo image2.png -> image1.png ---
$ftrFirst = $sec.Footers(2)
if ($ftrFirst.Exists -and $ftrFirst.Range.InlineShapes.Count -ge 1) {
    $pearsonLogoB = $ftrFirst.Range.InlineShapes(1)
    $pearsonLogoB.Name = "image1.png"
}
